$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "SUM of range" examples appended to the existing F/G label-formula table.
# Leading apostrophe forces text entry (matches the quotePrefix style already
# used by the sibling label cells F5/F6/F8) instead of being parsed as a formula.
$ws.Range("F9").Value = "'=SUM(B5:B8)"
$ws.Range("G9").Formula = "=SUM(B5:B8)"

$ws.Range("F10").Value = "'=SUM(B8,B9:B10)"
$ws.Range("G10").Formula = "=SUM(B8,B9:B10)"

# Defined names referencing empty anchor cells used while building the table
$wb.Names.Add("SUM", "=Sheet1!`$F`$23")
$wb.Names.Add("NOW", "=Sheet1!`$F`$24")

# Stray formatted-but-empty cells left over at the anchors used for the
# named ranges above; same date/time number format as B2 (NOW()).
$ws.Range("F24").NumberFormat = "m/d/yy h:mm"
$ws.Range("H30").NumberFormat = "m/d/yy h:mm"

# Column H picks up an explicit best-fit width from the new table entries
$ws.Columns.Item(8).ColumnWidth = 10.82

$ws.Range("F23:H31").Select()

$wb.Save()
